$d = $word.ActiveDocument
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- 1) "Enquete ..." -> "Enquête ..." (accent fix on the investigation bullet), typed
#        as two runs: "Enquête" then " sur de l'espionnage industriel qui tourne mal" ---
$enquetePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Enquete")) {
        $enquetePara = $p
        break
    }
}
$enqueteXml = ""
$enqueteXml += "<w:p $w w:rsidR='00AB50D8' w:rsidRDefault='001E0D80' w:rsidP='00AB50D8'>"
$enqueteXml +=   "<w:pPr>"
$enqueteXml +=     "<w:pStyle w:val='Paragraphedeliste'/>"
$enqueteXml +=     "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>"
$enqueteXml +=   "</w:pPr>"
$enqueteXml +=   "<w:r><w:t>Enquête</w:t></w:r>"
$enqueteXml +=   "<w:r><w:t xml:space='preserve'> sur de l’espionnage industriel qui tourne mal</w:t></w:r>"
$enqueteXml += "</w:p>"
$enquetePara.Range.InsertXML($enqueteXml)

# --- 2) Add the "Nom personnage principal" / "Usine" bullets after the last paragraph ---
# Locate the paragraph that ends the document ("...qui lui permettent ") which currently
# carries the (hidden) _GoBack bookmark at its very end.
$lastPara = $d.Paragraphs.Last

$xml = ""
$xml += "<w:p $w w:rsidR='003E39B0' w:rsidRPr='00E071C2' w:rsidRDefault='003E39B0' w:rsidP='000019F1'>"
$xml +=   "<w:pPr>"
$xml +=     "<w:pStyle w:val='Paragraphedeliste'/>"
$xml +=     "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>"
$xml +=   "</w:pPr>"
$xml +=   "<w:r><w:t xml:space='preserve'>En s’enfonçant dans l’usine, on récupère des gemmes sur les monstres qui lui permettent </w:t></w:r>"
$xml += "</w:p>"
$xml += "<w:p $w>"
$xml +=   "<w:pPr>"
$xml +=     "<w:pStyle w:val='Paragraphedeliste'/>"
$xml +=     "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>"
$xml +=     "<w:rPr><w:lang w:val='en-GB'/></w:rPr>"
$xml +=   "</w:pPr>"
$xml +=   "<w:r><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:t>Nom person</w:t></w:r>"
$xml +=   "<w:r><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:t>n</w:t></w:r>"
$xml +=   "<w:r><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:t>age principal:</w:t></w:r>"
$xml +=   "<w:r><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:t xml:space='preserve'> Wilfred Clifford Turner</w:t></w:r>"
$xml += "</w:p>"
$xml += "<w:p $w>"
$xml +=   "<w:pPr>"
$xml +=     "<w:pStyle w:val='Paragraphedeliste'/>"
$xml +=     "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>"
$xml +=     "<w:rPr><w:lang w:val='en-GB'/></w:rPr>"
$xml +=   "</w:pPr>"
$xml +=   "<w:r><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:t>Usine : Astrate Vapora</w:t></w:r>"
$xml +=   "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>"
$xml += "</w:p>"
$xml += "<w:p $w>"
$xml +=   "<w:pPr>"
$xml +=     "<w:pStyle w:val='Paragraphedeliste'/>"
$xml +=     "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>"
$xml +=     "<w:rPr><w:lang w:val='en-GB'/></w:rPr>"
$xml +=   "</w:pPr>"
$xml += "</w:p>"

# Replacing the whole last-paragraph range (rather than inserting at a collapsed point)
# drops the old hidden _GoBack bookmark that lived there and re-creates it only once,
# on the new "Usine" paragraph, matching the target layout.
$lastPara.Range.InsertXML($xml)

Write-Output "done"
